# Commit: "Fruta / hortaliza, semanal"
# Insert one new weekly record as row 101 in the data table, pushing the
# existing rows 101-124 down to 102-125 (the last former row, 124, ends up
# as the new row 125 - a plain append, no data loss).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 101 (shifts 101..124 -> 102..125,
# preserving cell styles/number formats such as the date style on column D).
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with this week's observation.
$ws.Range("A101").Value = 5
$ws.Range("B101").Value = "Macroferia Regional de Talca"
$ws.Range("C101").Value = "Maule"
$ws.Range("D101").Value = 44551
$ws.Range("E101").Value = 7
$ws.Range("F101").Value = 100112031
$ws.Range("G101").Value = "Poroto verde"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 300
$ws.Range("K101").Value = 23000
$ws.Range("L101").Value = 23000
$ws.Range("M101").Value = 23000
$ws.Range("N101").Value = "`$/saco 25 kilos"
$ws.Range("O101").Value = "Región del Maule"
$ws.Range("P101").Value = 920
$ws.Range("Q101").Value = 25
$ws.Range("R101").Value = "Hortaliza"
